$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.367497205734253
$ws.Range("B1").Value = 6.692091941833496
$ws.Range("C1").Value = 5.822466850280762
$ws.Range("D1").Value = 6.575465679168701
$ws.Range("E1").Value = 4.787083148956299
